$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17:93 down to 18:94
$ws.Rows("17:17").Insert()

# Fill in the new row 17 with the new observation.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are constant for every data row in this sheet,
# so copy them down from the row below (the row that used to be row 17).
$ws.Cells.Item(17, 1).Value  = $ws.Cells.Item(18, 1).Value2   # A: Mercado ID
$ws.Cells.Item(17, 2).Value  = $ws.Cells.Item(18, 2).Value2   # B: Mercado
$ws.Cells.Item(17, 3).Value  = $ws.Cells.Item(18, 3).Value2   # C: Región
$ws.Cells.Item(17, 4).Value  = 44623                          # D: Fecha
$ws.Cells.Item(17, 5).Value  = $ws.Cells.Item(18, 5).Value2   # E: Codreg
$ws.Cells.Item(17, 6).Value  = $ws.Cells.Item(18, 6).Value2   # F: Tipo
$ws.Cells.Item(17, 7).Value  = $ws.Cells.Item(18, 7).Value2   # G: Producto ID
$ws.Cells.Item(17, 8).Value  = $ws.Cells.Item(18, 8).Value2   # H: Producto
$ws.Cells.Item(17, 9).Value  = $ws.Cells.Item(18, 9).Value2   # I: Categoría ID
$ws.Cells.Item(17, 10).Value = $ws.Cells.Item(18, 10).Value2  # J: Categoría
$ws.Cells.Item(17, 11).Value = $ws.Cells.Item(18, 11).Value2  # K: Variedad
$ws.Cells.Item(17, 12).Value = "Primera"                      # L: Calidad
$ws.Cells.Item(17, 13).Value = 220                            # M: Volumen
$ws.Cells.Item(17, 14).Value = 2500                           # N: Precio mínimo
$ws.Cells.Item(17, 15).Value = 3000                           # O: Precio máximo
$ws.Cells.Item(17, 16).Value = 2773                           # P: Precio promedio ponderado
$ws.Cells.Item(17, 17).Value = $ws.Cells.Item(18, 17).Value2  # Q: Unidad de comercialización
$ws.Cells.Item(17, 18).Value = "Provincia de Linares"         # R: Origen
$ws.Cells.Item(17, 19).Value = 1386                           # S: Precio $/Kg
$ws.Cells.Item(17, 20).Value = $ws.Cells.Item(18, 20).Value2  # T: Kg / unidad

# The D column (Fecha) uses a custom date/time number format elsewhere in the
# column; make sure the newly inserted row keeps it (Insert() should already
# carry it from the row above, but set it explicitly to be safe).
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
